# Applies the "Actualizacion automatica" edit described by the diff:
#  - Sheet "VENTAS POR GRUPO": insert a new salesperson row ("GARCES MORALES
#    ANA CRISTINA") before "GARCIA BRAVO JOSE LUIS" (old row 288), plus a
#    handful of standalone value corrections.
#  - Sheet "VENTA MENSUAL": the mirror insertion (old row 292) plus the same
#    value corrections (transposed into its own layout) and updated totals.
#  - Sheet "CUMPLIMIENTO MENSUAL": recomputed VENTA / POR CUMPLIR /
#    CUMPLIMIENTO figures for a few asesor/grupo rows, plus a narrower
#    CUMPLIMIENTO column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Stand-alone value corrections (no row movement involved).
$ws1.Cells.Item(28, 5).Value = 207.38    # E28
$ws1.Cells.Item(74, 3).Value = 129.6     # C74
$ws1.Cells.Item(74, 13).Value = 811.75   # M74
$ws1.Cells.Item(107, 13).Value = -34.85  # M107

# Insert the new row for "GARCES MORALES ANA CRISTINA" right above the old
# row 288 ("GARCIA BRAVO JOSE LUIS"); this pushes every row below down by
# one (288-338 -> 289-339) and the old footer stats row (339 -> 340).
$ws1.Rows.Item(288).Insert()

$ws1.Cells.Item(288, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(288, 2).Value = "GARCES MORALES ANA CRISTINA"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(288, $col).Value = 0
}

# Refresh the "<n> de 337" -> "<n> de 338" footer counters (now row 340).
$ws1.Cells.Item(340, 3).Value = "3 de 338"
$ws1.Cells.Item(340, 4).Value = "4 de 338"
$ws1.Cells.Item(340, 5).Value = "1 de 338"
$ws1.Cells.Item(340, 6).Value = "0 de 338"
$ws1.Cells.Item(340, 7).Value = "0 de 338"
$ws1.Cells.Item(340, 8).Value = "3 de 338"
$ws1.Cells.Item(340, 9).Value = "4 de 338"
$ws1.Cells.Item(340, 10).Value = "0 de 338"
$ws1.Cells.Item(340, 11).Value = "1 de 338"
$ws1.Cells.Item(340, 12).Value = "4 de 338"
$ws1.Cells.Item(340, 13).Value = "10 de 338"
$ws1.Cells.Item(340, 14).Value = "0 de 338"
$ws1.Cells.Item(340, 15).Value = "0 de 338"
$ws1.Cells.Item(340, 16).Value = "1 de 338"
$ws1.Cells.Item(340, 17).Value = "0 de 338"
$ws1.Cells.Item(340, 18).Value = "0 de 338"

# ---------------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Stand-alone value corrections (same underlying figures as sheet 1, just
# laid out in this sheet's column scheme).
$ws2.Cells.Item(28, 6).Value = 207.38    # F28
$ws2.Cells.Item(74, 6).Value = 941.35    # F74
$ws2.Cells.Item(107, 6).Value = -34.85   # F107

# Same new-row insertion, mirrored at this sheet's old row 292 ("GARCIA
# BRAVO JOSE LUIS"), pushing 292-342 -> 293-343 and the totals row 343 -> 344.
$ws2.Rows.Item(292).Insert()

$ws2.Cells.Item(292, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(292, 2).Value = "GARCES MORALES ANA CRISTINA"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(292, $col).Value = 0
}

# Update the grand-total row (now row 344): only the F (POR CUMPLIR) column
# actually changes, the rest keep their previous totals.
$ws2.Cells.Item(344, 6).Value = 19739.5

# ---------------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ALMEIDA CUATIN JHONATHANN CARLOS / FREGADEROS DE COCINA
$ws3.Cells.Item(4, 4).Value = 207.38
$ws3.Cells.Item(4, 5).Value = 314.23144263264
$ws3.Cells.Item(4, 6).Value = 0.3975756339878713

# CASTRO ALCIVAR EDA MARIA / 240X120 PORCELANATO
$ws3.Cells.Item(14, 4).Value = 388.8
$ws3.Cells.Item(14, 5).Value = 2503.40588040374
$ws3.Cells.Item(14, 6).Value = 0.1344302639844315

# CASTRO ALCIVAR EDA MARIA / PORCELANATO
$ws3.Cells.Item(24, 4).Value = 767.97
$ws3.Cells.Item(24, 5).Value = 47856.09
$ws3.Cells.Item(24, 6).Value = 0.01579403283066038

# GUERRERO FAREZ FABIAN MAURICIO / PORCELANATO
$ws3.Cells.Item(36, 4).Value = 3848.15
$ws3.Cells.Item(36, 5).Value = 48814.97
$ws3.Cells.Item(36, 6).Value = 0.07307105997517807

# TOTAL row
$ws3.Cells.Item(76, 4).Value = 21032.64
$ws3.Cells.Item(76, 5).Value = 386579.3170193434
$ws3.Cells.Item(76, 6).Value = 0.05159966393969619

# Narrow the CUMPLIMIENTO column (F) from 28 to 25 raw character-width units.
# 24.12 is the ColumnWidth that this host's pixel-rounding maps onto a raw
# width of exactly 25.
$ws3.Columns.Item(6).ColumnWidth = 24.12
